# Add new pendulum/theta-controller experiments (rows 11-19) to the
# experiments schedule sheet: theta control with damped alpha, and the
# PP + LMI + int controller for theta and alpha.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Experiment 8 (row 11): high-angle pendulum swing, only the pendulum ---
$ws.Range("C11").Value = "High angles only pendulum swing"
$ws.Range("E11").Value = "Check non-linear characteristic of the pendulum"
$ws.Range("G11").Value = "Paolo"

# --- Experiment 9 (row 12): high-angle full pendulum swing ---
$ws.Range("C12").Value = "High angles full pendulum swing"

$ws.Range("D11").Value = "Connect the system, remove the pendulum, then move it to 90/180 deg angle and let it fall"

$ws.Range("D12").Value = "Same as 8 but attached to the cube"
$ws.Range("E12").Value = "Same as 8 for whole system"
$ws.Range("G12").Value = "Paolo"

# --- Experiment 10 (row 13): dynamic performance on theta of PP + LMI ---
$ws.Range("C13").Value = "Dynamic performance on theta of PP + LMI"
$ws.Range("D13").Value = "Use a sinesweep as input for theta on the CL system with full-state PP control"
$ws.Range("E13").Value = "Check dynamic performances of full state PP"
$ws.Range("G13").Value = "Paolo"

# --- Experiment 11 (row 14): static performance, theta controller + alpha damper ---
$ws.Range("C14").Value = "Static performance of theta controller + alpha damper"
$ws.Range("D14").Value = "Step input to CL system with theta + alpha damper controller"
$ws.Range("E14").Value = "Check static performances of theta controller + alpha damper"
$ws.Range("G14").Value = "Paolo"

# --- Experiment 12 (row 15): dynamic performance, theta controller + alpha damper ---
$ws.Range("C15").Value = "Dynamic performance of theta controller + alpha damper"
$ws.Range("D15").Value = "Sinesweep input to CL system with theta + alpha damper controller"
$ws.Range("E15").Value = "Check dynamic performances of theta controller + alpha damper"
$ws.Range("G15").Value = "Paolo"

# --- Experiment 13 (row 16): robustness, theta controller + alpha damper ---
$ws.Range("C16").Value = "Robustness of theta controller + alpha damper"
$ws.Range("D16").Value = "Step input to CL system with theta + alpha damper controller. Abuse when stable"
$ws.Range("E16").Value = "Check robustness of theta controller + alpha damper"
$ws.Range("G16").Value = "Paolo"

# --- Experiments 14-16 (rows 17-19): PP + LMI + int controller, filled by column ---
$ws.Range("C17").Value = "Static performance of PP + LMI + int"
$ws.Range("C18").Value = "Dynamic performance of PP + LMI + int"
$ws.Range("C19").Value = "Robustness of PP + LMI + int"

$ws.Range("D17").Value = "Step input to CL system with PP + LMI + int controller"
$ws.Range("D18").Value = "Sinsesweep input to CL system with PP + LMI + int controller"
$ws.Range("D19").Value = "Step input to CL system with PP + LMI + int controller. Abuse when stable"

$ws.Range("E17").Value = "Check static performances of PP + LMI + int"
$ws.Range("E18").Value = "Check dynamic performances of  PP + LMI + int"
$ws.Range("E19").Value = "Check robustness of  PP + LMI + int"

$ws.Range("G17").Value = "Paolo"
$ws.Range("G18").Value = "Paolo"
$ws.Range("G19").Value = "Paolo"
